# Update workbook with corrected forecast output:
# - Rename Sheet1 -> "Sales vs PO"
# - Add 3 new sheets: "Weekly Growth", "Volume Insights", "Prediction Info"
# - On "Sales vs PO": insert a new "Order Week" column holding the original
#   order-week dates, and shift the "ds" column forward to the forecast week end
# - Populate headers (and placeholder data) on the new sheets

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

# ---- Sales vs PO: insert "Order Week" column before PO_Requested_Qty ----
$ws1.Columns.Item(3).Insert()
$ws1.Cells.Item(1, 3).Value2 = "Order Week"
$ws1.Range("C2:C17").NumberFormat = "YYYY-MM-DD HH:MM:SS"

for ($r = 2; $r -le 17; $r++) {
    $origDs = $ws1.Cells.Item($r, 1).Value2
    $ws1.Cells.Item($r, 3).Value2 = $origDs
    $ws1.Cells.Item($r, 1).Value2 = $origDs + 6
}

# ---- Add the remaining sheets, in order, right after "Sales vs PO" ----
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"

$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"

# ---- Weekly Growth headers ----
$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)
$ws2.Cells.Item(1, 1).Value2 = "ds"
$ws2.Cells.Item(1, 2).Value2 = "PO_Requested_Qty"
$ws2.Cells.Item(1, 3).Value2 = "Growth%"

# ---- Volume Insights headers + placeholder data row ----
$ws1.Range("A1:D1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)
$ws3.Cells.Item(1, 1).Value2 = "Total_PO_Quantity"
$ws3.Cells.Item(1, 2).Value2 = "Average_PO_Quantity"
$ws3.Cells.Item(1, 3).Value2 = "Max_PO_Quantity"
$ws3.Cells.Item(1, 4).Value2 = "Min_PO_Quantity"

$ws3.Cells.Item(2, 1).Value2 = 0
$ws3.Cells.Item(2, 2).Value2 = 0
$ws3.Cells.Item(2, 3).Value2 = 0
$ws3.Cells.Item(2, 4).Value2 = 0

# ---- Prediction Info header + placeholder data row ----
$ws1.Range("A1:A1").Copy()
$ws4.Range("A1:A1").PasteSpecial(-4122)
$ws4.Cells.Item(1, 1).Value2 = "Predicted_Next_Week_PO_Quantity"

$ws4.Cells.Item(2, 1).Value2 = 0

# ---- Activate first sheet to mirror original workbook state ----
$ws1.Activate()
